$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 9, keeping only row 1
$ws.Range("A2:B9").EntireRow.Delete() | Out-Null

# Update row 1 values: the last fixture (Real Betis Balompié) replaces the first one,
# and its date is corrected from 22/05/2022 to 20/05/2022.
$ws.Range("A1").Value = "Real Madrid CF - Real Betis Balompié"
$ws.Range("B1").Value = "20/05/2022"

$wb.Save()
